$wb = $excel.ActiveWorkbook

$repay = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before the "Late" column (N) on the
# "Repayment schedule" sheet so a new column is available for data entry
# (Variable Instalments). This shifts old N->O, O->P, P->Q.
$repay.Columns("N:N").Insert() | Out-Null

# The newly inserted column should not be a "best fit" column like its
# neighbours; give it an explicit width (matches column M's width of 11).
$repay.Columns("N:N").ColumnWidth = 10.17

# Update the selection on the Repayment schedule sheet and make it the
# active sheet/tab.
$repay.Range("R6").Select() | Out-Null
$repay.Activate() | Out-Null
